$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update DAMSLTag (column I) and DialogAct (column J) values for the rows
# re-annotated by the SGNN dialog-act tagger re-run.
$updates = @(
    @{ Row = 4; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 5; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 16; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 29; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 42; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 61; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 65; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 68; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 80; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 93; DAMSLTag = 'ba'; DialogAct = 'Appreciation' }
    @{ Row = 109; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 117; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 120; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 125; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 135; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 141; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 142; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 147; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 149; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 153; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 154; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 155; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 156; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 164; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 170; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 177; DAMSLTag = 'ba'; DialogAct = 'Appreciation' }
    @{ Row = 178; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 181; DAMSLTag = 'qy'; DialogAct = 'Yes-No-Question' }
    @{ Row = 186; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 187; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 192; DAMSLTag = 'ba'; DialogAct = 'Appreciation' }
    @{ Row = 206; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 207; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 221; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 222; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 234; DAMSLTag = 'qy'; DialogAct = 'Yes-No-Question' }
    @{ Row = 240; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 246; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 247; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 252; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 274; DAMSLTag = 'ba'; DialogAct = 'Appreciation' }
    @{ Row = 276; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 283; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 284; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 285; DAMSLTag = 'ba'; DialogAct = 'Appreciation' }
    @{ Row = 295; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 297; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 305; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 314; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 342; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 351; DAMSLTag = '%'; DialogAct = 'Uninterpretable' }
    @{ Row = 353; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 355; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 357; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 360; DAMSLTag = 'qy'; DialogAct = 'Yes-No-Question' }
    @{ Row = 371; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 375; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 381; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 382; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 390; DAMSLTag = 'ba'; DialogAct = 'Appreciation' }
    @{ Row = 408; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 413; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 423; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 464; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 485; DAMSLTag = 'ba'; DialogAct = 'Appreciation' }
    @{ Row = 493; DAMSLTag = 'ba'; DialogAct = 'Appreciation' }
    @{ Row = 497; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 507; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 508; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 512; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 516; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 520; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 532; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 541; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 545; DAMSLTag = 'ba'; DialogAct = 'Appreciation' }
    @{ Row = 548; DAMSLTag = 'qy'; DialogAct = 'Yes-No-Question' }
    @{ Row = 551; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 554; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 556; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 562; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 571; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 577; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 580; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 611; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 616; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.DAMSLTag
    $ws.Cells.Item($u.Row, 10).Value = $u.DialogAct
}

$wb.Save()